$wb = $excel.ActiveWorkbook

foreach ($ws in @($wb.Worksheets.Item(1), $wb.Worksheets.Item(2))) {
    $ws.Range("F1").Value = "Mid Paper 1"
    $ws.Range("G1").Value = "Mid Paper 2"

    for ($row = 2; $row -le 4; $row++) {
        $ws.Cells.Item($row, 6).Value = $ws.Cells.Item($row, 4).Value2
        $ws.Cells.Item($row, 7).Value = $ws.Cells.Item($row, 5).Value2
    }
}

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws1.Range("F1:G4").Select()
$ws2.Activate()
$ws2.Range("F1:G4").Select()

Write-Host "done"
